$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il18"
$ws.Range("C2").Value = "Il18r1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.855689333333333
$ws.Range("H2").Value = 11.567068
$ws.Range("I2").Value = 0.1340401150840085
$ws.Range("J2").Value = 0.1340401150840085
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.024585
$ws.Range("N2").Value = 0.073755
$ws.Range("O2").Value = 0.08088989301323214
$ws.Range("P2").Value = 0.08088989301323216
$ws.Range("Q2").Value = 0.09479212225999999
$ws.Range("R2").Value = 0.8531291003399999
$ws.Range("S2").Value = 0.01084249056862677
$ws.Range("T2").Value = 0.01084249056862678

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il18"
$ws.Range("C3").Value = "Il18r1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.855689333333333
$ws.Range("H3").Value = 11.567068
$ws.Range("I3").Value = 0.1340401150840085
$ws.Range("J3").Value = 0.1340401150840085
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2793466666666666
$ws.Range("N3").Value = 0.8380399999999999
$ws.Range("O3").Value = 0.9191101069867678
$ws.Range("P3").Value = 0.9191101069867679
$ws.Range("Q3").Value = 1.077073962968889
$ws.Range("R3").Value = 9.693665666719998
$ws.Range("S3").Value = 0.1231976245153817
$ws.Range("T3").Value = 0.1231976245153818

# Row 4
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Il18"
$ws.Range("C4").Value = "Il18r1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 22.00607133333333
$ws.Range("H4").Value = 66.018214
$ws.Range("I4").Value = 0.7650243780187601
$ws.Range("J4").Value = 0.7650243780187601
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.024585
$ws.Range("N4").Value = 0.073755
$ws.Range("O4").Value = 0.08088989301323214
$ws.Range("P4").Value = 0.08088989301323216
$ws.Range("Q4").Value = 0.54101926373
$ws.Range("R4").Value = 4.86917337357
$ws.Range("S4").Value = 0.06188274009045197
$ws.Range("T4").Value = 0.06188274009045198

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Il18"
$ws.Range("C5").Value = "Il18r1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 22.00607133333333
$ws.Range("H5").Value = 66.018214
$ws.Range("I5").Value = 0.7650243780187601
$ws.Range("J5").Value = 0.7650243780187601
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2793466666666666
$ws.Range("N5").Value = 0.8380399999999999
$ws.Range("O5").Value = 0.9191101069867678
$ws.Range("P5").Value = 0.9191101069867679
$ws.Range("Q5").Value = 6.147322673395555
$ws.Range("R5").Value = 55.32590406055999
$ws.Range("S5").Value = 0.7031416379283081
$ws.Range("T5").Value = 0.7031416379283082

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Il18"
$ws.Range("C6").Value = "Il18r1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.903429
$ws.Range("H6").Value = 8.710287000000001
$ws.Range("I6").Value = 0.1009355068972313
$ws.Range("J6").Value = 0.1009355068972313
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.024585
$ws.Range("N6").Value = 0.073755
$ws.Range("O6").Value = 0.08088989301323214
$ws.Range("P6").Value = 0.08088989301323216
$ws.Range("Q6").Value = 0.071380801965
$ws.Range("R6").Value = 0.6424272176850001
$ws.Range("S6").Value = 0.008164662354153395
$ws.Range("T6").Value = 0.008164662354153397

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Il18"
$ws.Range("C7").Value = "Il18r1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.903429
$ws.Range("H7").Value = 8.710287000000001
$ws.Range("I7").Value = 0.1009355068972313
$ws.Range("J7").Value = 0.1009355068972313
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2793466666666666
$ws.Range("N7").Value = 0.8380399999999999
$ws.Range("O7").Value = 0.9191101069867678
$ws.Range("P7").Value = 0.9191101069867679
$ws.Range("Q7").Value = 0.8110632130533334
$ws.Range("R7").Value = 7.29956891748
$ws.Range("S7").Value = 0.09277084454307791
$ws.Range("T7").Value = 0.09277084454307792
